# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" for the 87efa192 file row on the
# per-locale status sheets (row 4, column D) to reflect a new handoff run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-23 08:33:39"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-23 08:33:50"
